$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data describing the two new localization entries that are being
# inserted ahead of the existing "d70603bd..." row on every sheet.
# ---------------------------------------------------------------------------

$newDateOverview = "2016-17-13 12:17:55"
$newDateZh = "2016-03-13 12:17:51"
$newDateDe = "2016-03-13 12:17:55"

# =============================================================================
# Sheet "Overview"
# =============================================================================
$ws = $wb.Worksheets.Item("Overview")

# Make room for the two new rows (700ca9e0..., d2f5ec2f...) right before the
# existing "d70603bd..." row (currently row 3).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Refresh the handoff date text on the pre-existing rows.
$ws.Range("D2").Value = $newDateOverview
$ws.Range("D5").Value = $newDateOverview

# New row for 700ca9e0...
$ws.Range("A3").Value = "700ca9e0-6d9c-4b79-9354-073a22096042.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $newDateOverview

# New row for d2f5ec2f...
$ws.Range("A4").Value = "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = $newDateOverview

# Hyperlinks shift around because of the inserted rows - rebuild them in the
# correct final order.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/35fd21a0-345a-48db-befc-ccf840d57e50.md", "", "", "35fd21a0-345a-48db-befc-ccf840d57e50.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/700ca9e0-6d9c-4b79-9354-073a22096042.md", "", "", "700ca9e0-6d9c-4b79-9354-073a22096042.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md", "", "", "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d70603bd-12d3-4f8e-aa0b-af284c8244fe.md", "", "", "d70603bd-12d3-4f8e-aa0b-af284c8244fe.md")

# =============================================================================
# Sheet "zh-cn"
# =============================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Refresh handoff datetime text for the pre-existing rows.
$ws.Range("E2").Value = $newDateZh
$ws.Range("E5").Value = $newDateZh

# New row for 700ca9e0...
$ws.Range("A3").Value = "700ca9e0-6d9c-4b79-9354-073a22096042.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "700ca9e0-6d9c-4b79-9354-073a22096042.8f8a606ae8f704f6e0abe792b2761c0362d76970.zh-cn.xlf"
$ws.Range("E3").Value = $newDateZh
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# New row for d2f5ec2f...
$ws.Range("A4").Value = "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.ea765e34ab8675f792c0b51b4913986e98ac9713.zh-cn.xlf"
$ws.Range("E4").Value = $newDateZh
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("I4").Value = "Include"

# Hyperlinks shift around because of the inserted rows - rebuild them in the
# correct final order.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/35fd21a0-345a-48db-befc-ccf840d57e50.md", "", "", "35fd21a0-345a-48db-befc-ccf840d57e50.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/35fd21a0-345a-48db-befc-ccf840d57e50.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3b342bbdeb217743d072a5572888f79e5df662f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/35fd21a0-345a-48db-befc-ccf840d57e50.9870b55b67b629256b6c9329b32ca6cd0a25df80.zh-cn.xlf", "", "", "35fd21a0-345a-48db-befc-ccf840d57e50.9870b55b67b629256b6c9329b32ca6cd0a25df80.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/700ca9e0-6d9c-4b79-9354-073a22096042.md", "", "", "700ca9e0-6d9c-4b79-9354-073a22096042.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/700ca9e0-6d9c-4b79-9354-073a22096042.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3b342bbdeb217743d072a5572888f79e5df662f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/700ca9e0-6d9c-4b79-9354-073a22096042.8f8a606ae8f704f6e0abe792b2761c0362d76970.zh-cn.xlf", "", "", "700ca9e0-6d9c-4b79-9354-073a22096042.8f8a606ae8f704f6e0abe792b2761c0362d76970.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md", "", "", "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3b342bbdeb217743d072a5572888f79e5df662f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.ea765e34ab8675f792c0b51b4913986e98ac9713.zh-cn.xlf", "", "", "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.ea765e34ab8675f792c0b51b4913986e98ac9713.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d70603bd-12d3-4f8e-aa0b-af284c8244fe.md", "", "", "d70603bd-12d3-4f8e-aa0b-af284c8244fe.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d70603bd-12d3-4f8e-aa0b-af284c8244fe.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3b342bbdeb217743d072a5572888f79e5df662f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d70603bd-12d3-4f8e-aa0b-af284c8244fe.b87f5d5568c96865137e31271d1d42e8984ef81c.zh-cn.xlf", "", "", "d70603bd-12d3-4f8e-aa0b-af284c8244fe.b87f5d5568c96865137e31271d1d42e8984ef81c.zh-cn.xlf")

# =============================================================================
# Sheet "de-de"
# =============================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Refresh handoff datetime text for the pre-existing rows.
$ws.Range("E2").Value = $newDateDe
$ws.Range("E5").Value = $newDateDe

# New row for 700ca9e0...
$ws.Range("A3").Value = "700ca9e0-6d9c-4b79-9354-073a22096042.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "700ca9e0-6d9c-4b79-9354-073a22096042.8f8a606ae8f704f6e0abe792b2761c0362d76970.de-de.xlf"
$ws.Range("E3").Value = $newDateDe
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# New row for d2f5ec2f...
$ws.Range("A4").Value = "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.ea765e34ab8675f792c0b51b4913986e98ac9713.de-de.xlf"
$ws.Range("E4").Value = $newDateDe
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("I4").Value = "Include"

# Hyperlinks shift around because of the inserted rows - rebuild them in the
# correct final order.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/35fd21a0-345a-48db-befc-ccf840d57e50.md", "", "", "35fd21a0-345a-48db-befc-ccf840d57e50.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/35fd21a0-345a-48db-befc-ccf840d57e50.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e36d1593f1a8795e26c44f52aa39efe380285b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/35fd21a0-345a-48db-befc-ccf840d57e50.9870b55b67b629256b6c9329b32ca6cd0a25df80.de-de.xlf", "", "", "35fd21a0-345a-48db-befc-ccf840d57e50.9870b55b67b629256b6c9329b32ca6cd0a25df80.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/700ca9e0-6d9c-4b79-9354-073a22096042.md", "", "", "700ca9e0-6d9c-4b79-9354-073a22096042.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/700ca9e0-6d9c-4b79-9354-073a22096042.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e36d1593f1a8795e26c44f52aa39efe380285b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/700ca9e0-6d9c-4b79-9354-073a22096042.8f8a606ae8f704f6e0abe792b2761c0362d76970.de-de.xlf", "", "", "700ca9e0-6d9c-4b79-9354-073a22096042.8f8a606ae8f704f6e0abe792b2761c0362d76970.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md", "", "", "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e36d1593f1a8795e26c44f52aa39efe380285b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.ea765e34ab8675f792c0b51b4913986e98ac9713.de-de.xlf", "", "", "d2f5ec2f-9a3c-48de-8d8c-8ebd3c1c202c.ea765e34ab8675f792c0b51b4913986e98ac9713.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d70603bd-12d3-4f8e-aa0b-af284c8244fe.md", "", "", "d70603bd-12d3-4f8e-aa0b-af284c8244fe.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/60925f2b1a046eba96b8f3a3e9b1abe4eaee7da5/e2e/d70603bd-12d3-4f8e-aa0b-af284c8244fe.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e36d1593f1a8795e26c44f52aa39efe380285b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d70603bd-12d3-4f8e-aa0b-af284c8244fe.b87f5d5568c96865137e31271d1d42e8984ef81c.de-de.xlf", "", "", "d70603bd-12d3-4f8e-aa0b-af284c8244fe.b87f5d5568c96865137e31271d1d42e8984ef81c.de-de.xlf")
